$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to remain Text (matches source data: textual prices)
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns per latest snapshot
$ws.Range("D2").Value = "28.601.14"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.880.67"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -3.15%  "
$ws.Range("D5").Value = "315.06"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "0.08435"
$ws.Range("D10").Value = "1.113"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").Value = "41.66"
$ws.Range("D12").Value = "6.269"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "1.875.17"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "7.275"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").Value = "0.00001106"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "91.21"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "0.06738"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").Value = "5.967"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").Value = "28.615.71"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "2.086.88"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").Value = "161.05"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "20.83"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").Value = "2.385"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").Value = "127.32"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "0.1052"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "5.797"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").Value = "3.608"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "0.02457"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").Value = "0.06526"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").Value = "0.2193"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").Value = "8.932"
$ws.Range("E38").Value = "  -5.85%  "
$ws.Range("D39").Value = "1.263"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "5.084"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "0.6447"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").Value = "11.19"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "0.6071"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "13.07"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "3.695"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "122.63"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "1.211"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  -8.90%  "
